$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new checklist item as the next row after the existing data.
$ws.Range("A21").Value = "Guid in ef core"

# Mirror Excel's natural post-entry selection: cursor moves to the next
# empty cell below what was just typed.
$ws.Range("A22").Select()
